# Adapt column header formatting to respective input file names (#7)
# - Rename header row (row 1) cells from *_old / *_new to *_FV2304 / *_FV2310
# - Turn the data range into a native Excel Table ("Table1") with an
#   autofilter and banded rows
# - Freeze the header row (top row) in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

$fv2310Headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

# Columns A-J hold the "_old" -> "_FV2304" headers, columns L-U hold the
# "_new" -> "_FV2310" headers (column K is the unchanged "diff" column).
for ($i = 0; $i -lt $fv2304Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2304Headers[$i]
}

for ($i = 0; $i -lt $fv2310Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2310Headers[$i]
}

# Turn the used data range into a real Excel table with an autofilter and
# banded rows ("Table1" over A1:U81).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$lastCol = 21
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))

$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze the header row.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
